# Update the date shown on the first slide from "24 July 2025" to
# "21 July 2025".
#
# The shape (id=100, "24 July 2025") holds the date as a single run.
# We only retype the day-of-month portion ("24 " -> "21 "), which is
# exactly what happens when a human edits the existing text in place in
# PowerPoint: the run gets split into the edited prefix and the
# untouched remainder ("July 2025").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Sanity check - make sure we are editing the expected shape/text.
if ($tr.Text -eq "24 July 2025") {
    # Replace just the leading "24 " (day-of-month + following space)
    # with "21 ", leaving "July 2025" as the untouched remainder so the
    # run splits exactly like a manual retype would.
    $dayPart = $tr.Characters(1, 3)
    $dayPart.Text = "21 "
}
